$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The original column A (numeric id values 3,5,7,13 with header-style border)
# is removed entirely, and every other column (B:F) shifts one place to the
# left (becoming A:E). Deleting column A reproduces that shift exactly.
$ws.Range("A:A").Delete()
